$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.647.90'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '1.642.80'
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.505'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.29%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.01'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +0.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0627'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.26'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.37%  '
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('E12').Value = '  +0.58%  '
$ws.Range('D13').Value = '1.647.49'
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('E14').Value = '  +2.11%  '
$ws.Range('E15').Value = '  +1.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.42'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.91%  '
$ws.Range('D17').Value = '26.694.26'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '217.03'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.35'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.30'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.27%  '
$ws.Range('E23').Value = '  +1.59%  '
$ws.Range('E24').Value = '  +14.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.53'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('E28').Value = '  +4.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.76'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.39%  '
$ws.Range('E30').Value = '  +2.52%  '
$ws.Range('E31').Value = '  +0.74%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.38'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.04'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.21%  '
$ws.Range('D34').Value = '1.278.21'
$ws.Range('E34').Value = '  +4.16%  '
$ws.Range('E35').Value = '  +2.78%  '
$ws.Range('E36').Value = '  +4.92%  '
$ws.Range('E37').Value = '  +0.47%  '
$ws.Range('E38').Value = '  +6.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.828'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.90%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.01'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('E41').Value = '  +2.50%  '
$ws.Range('E42').Value = '  -1.40%  '
$ws.Range('E43').Value = '  +2.22%  '
$ws.Range('D44').Value = '1.782.27'
$ws.Range('E44').Value = '  +0.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '91.85'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.74'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +8.11%  '
$ws.Range('E47').Value = '  +1.66%  '
$ws.Range('E48').Value = '  +0.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.82'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0970'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.08%  '
$ws.Range('E51').Value = '  -0.81%  '
